$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1872
$ws.Range("F5").Value = 70
$ws.Range("F6").Value = 718
$ws.Range("F7").Value = 102
$ws.Range("F8").Value = 493
$ws.Range("F9").Value = 853
$ws.Range("F10").Value = 1573
$ws.Range("F11").Value = 1257
$ws.Range("F12").Value = 1499
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 1450
$ws.Range("F15").Value = 334
$ws.Range("F16").Value = 1671
$ws.Range("F17").Value = 800
$ws.Range("F18").Value = 1083
$ws.Range("F19").Value = 358
$ws.Range("F21").Value = 109
$ws.Range("F22").Value = 1650
$ws.Range("F23").Value = 203
$ws.Range("F24").Value = 813
$ws.Range("F25").Value = 555
$ws.Range("F26").Value = 1180
$ws.Range("F27").Value = 302395
$ws.Range("F28").Value = 1033
$ws.Range("F29").Value = 64
$ws.Range("F30").Value = 568
$ws.Range("F32").Value = 1127
$ws.Range("F33").Value = 900
$ws.Range("F35").Value = 1118
$ws.Range("F36").Value = 66
$ws.Range("F37").Value = 247
$ws.Range("F39").Value = 869
$ws.Range("F40").Value = 1664
$ws.Range("F41").Value = 4
$ws.Range("F42").Value = 112
$ws.Range("F43").Value = 6
$ws.Range("F44").Value = 80
$ws.Range("F45").Value = 822
$ws.Range("F47").Value = 794
$ws.Range("F48").Value = 116
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 46
$ws.Range("G5").Value = 80
$ws.Range("F6").Value = 1490
$ws.Range("F7").Value = 75
$ws.Range("F9").Value = 2571
$ws.Range("F10").Value = 1209
$ws.Range("F11").Value = 408
$ws.Range("F12").Value = 723
$ws.Range("F14").Value = 31
$ws.Range("F18").Value = 452
$ws.Range("F22").Value = 81095
$ws.Range("F23").Value = 22
$ws.Range("F24").Value = 1
$ws.Range("F30").Value = 206
$ws.Range("F42").Value = 135
$ws.Range("F43").Value = 62
$ws.Range("F44").Value = 4
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 117
$ws.Range("F4").Value = 260
$ws.Range("F5").Value = 2873
$ws.Range("F6").Value = 4626
$ws.Range("F7").Value = 132
$ws.Range("F9").Value = 569
$ws.Range("F10").Value = 718
$ws.Range("F11").Value = 457
$ws.Range("F12").Value = 328
$ws.Range("F13").Value = 1031
$ws.Range("F14").Value = 273
$ws.Range("F15").Value = 636
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1872
$ws.Range("F3").Value = 260
$ws.Range("F4").Value = 4626
$ws.Range("F5").Value = 718
$ws.Range("F6").Value = 70
$ws.Range("F7").Value = 328
$ws.Range("F8").Value = 328
$ws.Range("F9").Value = 1031
$ws.Range("F10").Value = 1031
$ws.Range("F11").Value = 493
$ws.Range("F12").Value = 853
$ws.Range("F13").Value = 2571
$ws.Range("F14").Value = 1209
$ws.Range("F15").Value = 1573
$ws.Range("F16").Value = 1257
$ws.Range("F17").Value = 1499
$ws.Range("F18").Value = 50
$ws.Range("F19").Value = 1450
$ws.Range("F21").Value = 334
$ws.Range("F23").Value = 1671
$ws.Range("F24").Value = 800
$ws.Range("F25").Value = 1083
$ws.Range("F26").Value = 358
$ws.Range("F27").Value = 636
$ws.Range("F28").Value = 636
$ws.Range("F29").Value = 452
$ws.Range("F30").Value = 1650
$ws.Range("F32").Value = 203
$ws.Range("F33").Value = 813
$ws.Range("F34").Value = 555
$ws.Range("F35").Value = 1180
$ws.Range("F37").Value = 1033
$ws.Range("F38").Value = 64
$ws.Range("F39").Value = 568
$ws.Range("F40").Value = 1127
$ws.Range("F41").Value = 900
$ws.Range("F42").Value = 1118
$ws.Range("F44").Value = 247
$ws.Range("F45").Value = 869
$ws.Range("F47").Value = 1664
$ws.Range("F48").Value = 112
$ws.Range("F49").Value = 80
$ws.Range("F50").Value = 822
$ws.Range("F52").Value = 794
